{"js": "// 1. Rewrite the \"Heading 5\" paragraph text as two runs: \"Heading \" + \"5\"\n//    (drops the stray lastRenderedPageBreak that lived in the old single run\n//     and shortens \"Heading 5 - this is for the page break with MS Word\"\n//     down to just \"Heading 5\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/style\");\nawait context.sync();\n\nlet headingPara = null;\nfor (const p of paragraphs.items) {\n  if (p.style === \"Heading 5\") {\n    headingPara = p;\n    break;\n  }\n}\n\nconst headingRange = headingPara.getRange();\nheadingRange.insertText(\"Heading \", \"Replace\");\nawait context.sync();\n// Append \"5\" right after the range we just replaced, producing a second run.\nheadingRange.insertText(\"5\", \"After\");\nawait context.sync();\n\n// 2. Update the \"Heading 5\" style definition.\nconst style = context.document.getStyles().getByNameOrNullObject(\"Heading 5\");\nstyle.load(\"nameLocal\");\nawait context.sync();\n\n// Remove the forced page break before the heading.\nstyle.paragraphFormat.pageBreakBefore = false;\n// Change the heading font color from white/Background1 to black/Text1.\nstyle.font.color = \"#000000\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- 1. Rewrite the \"Heading 5\" paragraph text as two runs: \"Heading \" + \"5\" ---\n# (drops the stray lastRenderedPageBreak that lived in the old single run and\n#  shortens \"Heading 5 - this is for the page break with MS Word\" to \"Heading 5\")\n$headingPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Heading 5\") {\n        $headingPara = $p\n        break\n    }\n}\n\n$r = $headingPara.Range\n# Exclude the trailing paragraph mark from the range we rewrite.\n$r.End = $r.End - 1\n$r.Text = \"Heading \"\n# Append \"5\" as its own run right after the text we just set, so the\n# paragraph ends up as two separate runs instead of being merged into one.\n$r2 = $d.Range($r.End, $r.End)\n$r2.InsertAfter(\"5\")\n\n# --- 2. Update the \"Heading 5\" style definition ---\n$style = $d.Styles(\"Heading 5\")\n# Remove the forced page break before the heading.\n$style.ParagraphFormat.PageBreakBefore = 0\n# Change the heading font color from white/Background1 to black/Text1.\n$style.Font.TextColor.ObjectThemeColor = 13\n"}
